# Pos Table Management almost finished
# Updates the "Table List" sheet: reorders/extends the "Ground Floor" area
# rows and adds a full set of "AutoEdited" area rows (T1-T10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table List")

# --- Grow the "Ground Floor" block from 4 rows to 5 rows -------------------
# Existing rows 3:6 (Ground Floor / S3,S2,xyz,S1) shift down to 4:7.
$ws.Rows(3).Insert()

# --- Insert 6 new rows right before the existing "AutoEdited / T7" row -----
# That row currently sits at row 10 (after the insert above); after this
# insert it lands at row 16, matching the final layout.
$ws.Rows("10:15").Insert()

# --- (Re)write the "Ground Floor" area rows (3:7) in their final order -----
$groundFloor = @(
  @(3, "G1", 44),
  @(4, "S1", 11),
  @(5, "S3", 33),
  @(6, "S4", 55),
  @(7, "xyz", 77)
)
foreach ($item in $groundFloor) {
  $r = $item[0]
  $ws.Cells.Item($r, 1).Value = "Ground Floor"
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Cells.Item($r, 3).Value = $item[2]
  $ws.Cells.Item($r, 3).HorizontalAlignment = -4131
}

# --- New "AutoEdited" area rows (10:15 and 17:19), plus the pre-existing ---
# "AutoEdited / T7" row (now at row 16) whose Seating Capacity becomes a
# real left-aligned number instead of the old text value "47".
$autoEdited = @(
  @(10, "T1", 11),
  @(11, "T2", 22),
  @(12, "T3", 23),
  @(13, "T4", 44),
  @(14, "T5", 55),
  @(15, "T6", 66),
  @(16, "T7", 47),
  @(17, "T8", 88),
  @(18, "T9", 99),
  @(19, "T10", 10)
)
foreach ($item in $autoEdited) {
  $r = $item[0]
  $ws.Cells.Item($r, 1).Value = "AutoEdited"
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Cells.Item($r, 3).Value = $item[2]
  $ws.Cells.Item($r, 3).HorizontalAlignment = -4131
}

# --- Sheet view tidy-up: drop the frozen/scrolled top-left cell and move the
# selection to D17, matching the final saved view state.
$ws.Activate() | Out-Null
$ws.Range("D17").Select() | Out-Null
